{"js": "// The first two multiple-choice questions in this quiz-answer document\n// have their option runs prefixed with \"A. \", \"B. \", \"C. \", \"D. \".\n// This edit strips just that leading \"X. \" label (3 characters) from\n// the start of each of those 8 option paragraphs, leaving the answer\n// text, its run formatting (bold/highlight/font/etc.) and every other\n// run in the paragraph (tabs, trailing runs, ...) untouched.\n//\n// Later questions in the document reuse the exact same option letters\n// (\"A. \", \"B. \", ...) and some of the same words (e.g. \"Li\u00ean X\u00f4\"), so\n// matching is anchored to a short, distinctive substring from each of\n// the two target questions' option text and scoped to a single\n// paragraph before any text is removed - this keeps the edit limited\n// to exactly the 8 runs touched by the source change.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Distinctive substrings that uniquely identify each of the 8 target\n// option paragraphs (in document order).\nconst markers = [\n  \"ho\u00e0n to\u00e0n k\u1ebft th\u00fac\",\n  \"b\u01b0\u1edbc v\u00e0o giai \u0111o\u1ea1n k\u1ebft th\u00fac\",\n  \"\u0111ang di\u1ec5n ra v\u00f4 c\u00f9ng \u00e1c li\u1ec7t\",\n  \"b\u00f9ng n\u1ed5 v\u00e0 ng\u00e0y c\u00e0ng lan r\u1ed9ng\",\n  \"Anh, Ph\u00e1p, M\u1ef9\",\n  \"\u0110\u1ee9c, Italia, Nh\u1eadt\",\n  \"Anh, Ph\u00e1p, Li\u00ean X\u00f4\",\n  \"Li\u00ean X\u00f4, M\u1ef9, Anh\",\n];\n\nfor (const marker of markers) {\n  const paragraph = paragraphs.items.find((p) => p.text.indexOf(marker) !== -1);\n  if (!paragraph) {\n    throw new Error(\"Could not locate option paragraph for marker: \" + marker);\n  }\n\n  // Only the leading \"A. \" / \"B. \" / \"C. \" / \"D. \" (a letter, a period,\n  // a space - 3 characters) needs to go; everything else in the\n  // paragraph (answer text, tabs, later runs) stays exactly as-is.\n  const prefix = paragraph.text.substring(0, 3);\n  if (!/^[A-D]\\.\\s$/.test(prefix)) {\n    // Already edited (or unexpected shape) - skip rather than corrupt it.\n    continue;\n  }\n\n  const hits = paragraph.search(prefix, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly one prefix match in paragraph for marker '\" +\n        marker +\n        \"', found \" +\n        hits.items.length\n    );\n  }\n\n  hits.items[0].insertText(\"\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The first two multiple-choice questions in this quiz-answer document\n# have their option runs prefixed with \"A. \", \"B. \", \"C. \", \"D. \".\n# This edit strips just that leading \"X. \" label (3 characters) from\n# the start of each of those 8 option paragraphs, leaving the answer\n# text, its run formatting (bold/highlight/font/etc.) and every other\n# run in the paragraph (tabs, trailing runs, ...) untouched.\n#\n# Later questions in the document reuse the exact same option letters\n# (\"A. \", \"B. \", ...) and some of the same words (e.g. \"Li\u00ean X\u00f4\"), so\n# each target paragraph is first located via a short, distinctive\n# substring from its answer text, and the \"X. \" removal is done with\n# Find.Execute scoped to that single paragraph's Range - this keeps the\n# edit limited to exactly the 8 runs touched by the source change and\n# immune to the duplicate \"D. Li\u00ean X\u00f4, \" text later in the document.\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n\n# Distinctive substrings that uniquely identify each of the 8 target\n# option paragraphs (in document order).\n$markers = @(\n    \"ho\u00e0n to\u00e0n k\u1ebft th\u00fac\",\n    \"b\u01b0\u1edbc v\u00e0o giai \u0111o\u1ea1n k\u1ebft th\u00fac\",\n    \"\u0111ang di\u1ec5n ra v\u00f4 c\u00f9ng \u00e1c li\u1ec7t\",\n    \"b\u00f9ng n\u1ed5 v\u00e0 ng\u00e0y c\u00e0ng lan r\u1ed9ng\",\n    \"Anh, Ph\u00e1p, M\u1ef9\",\n    \"\u0110\u1ee9c, Italia, Nh\u1eadt\",\n    \"Anh, Ph\u00e1p, Li\u00ean X\u00f4\",\n    \"Li\u00ean X\u00f4, M\u1ef9, Anh\"\n)\n\nforeach ($marker in $markers) {\n    $target = $null\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $para = $d.Paragraphs($i)\n        if ($para.Range.Text.Contains($marker)) {\n            $target = $para\n            break\n        }\n    }\n\n    if ($null -eq $target) {\n        throw \"Could not locate option paragraph for marker: $marker\"\n    }\n\n    # Only the leading \"A. \" / \"B. \" / \"C. \" / \"D. \" (a letter, a\n    # period, a space - 3 characters) needs to go; everything else in\n    # the paragraph (answer text, tabs, later runs) stays exactly\n    # as-is.\n    $prefix = $target.Range.Text.Substring(0, 3)\n    if ($prefix -notmatch \"^[A-D]\\. $\") {\n        # Already edited (or unexpected shape) - skip rather than\n        # corrupt it.\n        continue\n    }\n\n    $rng = $target.Range\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($prefix, $false, $false, $false, $false, $false, $true, 1, $false, \"\", $wdReplaceOne)\n}\n"}
